$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text is a plain decimal number; briefly switch them to
# Text format so Excel stores the literal string (matching the scraped source
# data, e.g. "144.00") instead of silently normalising it to a numeric value,
# then restore the default style so formatting stays untouched.
$textCells = @("D5", "D6", "D8", "D11", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D31", "D33", "D37", "D38", "D39", "D40", "D43", "D44", "D45", "D46", "D48", "D51")
foreach ($cellRef in $textCells) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
}

$ws.Range("D2").Value = "61.469.90"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "2.991.28"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "599.41"
$ws.Range("E5").Value = "  +3.08%  "
$ws.Range("D6").Value = "144.00"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "2.989.77"
$ws.Range("E9").Value = "  -0.27%  "
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").Value = "6.06"
$ws.Range("E11").Value = "  +7.65%  "
$ws.Range("E12").Value = "  +3.87%  "
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "34.33"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "3.491.35"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "6.96"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "61.478.66"
$ws.Range("E18").Value = "  -1.17%  "
$ws.Range("D19").Value = "2.992.80"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").Value = "450.30"
$ws.Range("E20").Value = "  -0.70%  "
$ws.Range("D21").Value = "13.99"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "0.685"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").Value = "7.32"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "81.50"
$ws.Range("D25").Value = "10.81"
$ws.Range("E25").Value = "  +7.14%  "
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  -3.08%  "
$ws.Range("D27").Value = "11.95"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").Value = "2.68"
$ws.Range("E29").Value = "  +3.08%  "
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "7.19"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").Value = "27.29"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("E34").Value = "  +3.32%  "
$ws.Range("D35").Value = "0.0₃0829"
$ws.Range("E35").Value = "  +4.81%  "
$ws.Range("E36").Value = "  -0.52%  "
$ws.Range("D37").Value = "5.78"
$ws.Range("E37").Value = "  +1.24%  "
$ws.Range("D38").Value = "9.20"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("D39").Value = "50.39"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("D40").Value = "2.06"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("E41").Value = "  +11.15%  "
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "398.23"
$ws.Range("E43").Value = "  -3.79%  "
$ws.Range("D44").Value = "39.94"
$ws.Range("E44").Value = "  +5.16%  "
$ws.Range("D45").Value = "0.0353"
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "0.270"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").Value = "2.689.21"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("D48").Value = "130.95"
$ws.Range("E48").Value = "  +1.92%  "
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("E50").Value = "  -0.37%  "
$ws.Range("D51").Value = "2.14"
$ws.Range("E51").Value = "  +1.82%  "

# Restore the original (default) style on the cells we reformatted above.
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).Style = "Normal"
}
